$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.4
$ws.Range("F3").Value = 0.2
$ws.Range("F4").Value = 0.4
$ws.Range("F5").Value = 0

[void]$ws.Range("G15").Select()
